$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: H1 heading + bolded run near the end)
Replace-All "Play Fa Cai Shen Deluxe Slot for Free - Read Our Review" "Play Fa Cai Shen Deluxe Free | Slot Game Review"

# "What we like" bullets
Replace-All "Expanding Wild increases winning chances" "Suitable for all players"
Replace-All "Free Spins feature with random symbol transformation" "Asian-themed symbols and visuals"
Replace-All "Asian theme and matching soundtrack immerse players" "Expanding Wild symbol"
Replace-All "Max payout of 800x total stake" "Generous free spins feature"

# "What we don't like" bullets
Replace-All "High volatility may not suit all players" "High volatility gameplay"
Replace-All "Limited bonus features compared to other slot games" "Limited bet range"

# Meta description (italic run)
Replace-All "Read our review of Fa Cai Shen Deluxe online slot game. Play for free and discover the game's features, theme, symbols, and jackpot potential." "Play Fa Cai Shen Deluxe for free and read our review of this Asian-themed slot game."
